$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet (internal sheet name -> tab name)
$ws.Name = "alpha1F"

# Update the slightly recalculated values in row 13
$ws.Range("C13").Value = 0.9894216124101333
$ws.Range("D13").Value = 0.9939417127373235
$ws.Range("E13").Value = 0.9907191400909077
$ws.Range("F13").Value = 0.9894216124101333
$ws.Range("J13").Value = 0.9939417127373235
$ws.Range("K13").Value = 0.9923304264141155
$ws.Range("L13").Value = 0.9908760194121244
$ws.Range("M13").Value = 0.992431880614394
